$d = $word.ActiveDocument

# The paragraph currently reads "Version 2." built from the runs:
#   "Versi" | "on" | " 2" | <bookmark _GoBack> | "."
# It needs to become "Version 1." built from the runs:
#   "Version" | " 1." | <bookmark _GoBack>
# (i.e. revert commit "Wireframes version 2.")

# Step 1: the run holding " 2" (right before the bookmark) becomes " 1."
$rngSpace2 = $d.Content
$rngSpace2.Find.ClearFormatting()
$rngSpace2.Find.Execute(" 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngSpace2.Text = " 1."

# Step 2: merge the two runs spelling "Version" ("Versi" + "on") back into one run.
# Remove the "on" run, then rewrite the "Versi" run's text to "Version".
$rngOn = $d.Content
$rngOn.Find.ClearFormatting()
$rngOn.Find.Execute("on", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngOn.Delete()

$rngVersi = $d.Content
$rngVersi.Find.ClearFormatting()
$rngVersi.Find.Execute("Versi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngVersi.Text = "Version"

# Step 3: drop the now-orphaned "." run that used to trail the bookmark
# (text currently reads "Version 1.." - remove the stray final period character).
$endPos = $d.Content.End
$rngTrailingDot = $d.Range($endPos - 2, $endPos - 1)
$rngTrailingDot.Delete()
